# Fixing and re-run SCD0338-001 - SCD0338-014
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# P3: "10:00 PM" -> "10:30 PM" (leading apostrophe keeps the cell's
# existing text/quote-prefix formatting instead of Excel re-guessing it)
$ws.Range("P3").Value = "'10:30 PM"

# L2 / L3: "Test Daily Activity 3" -> "Test Daily Activity 5"
$ws.Range("L2").Value = "Test Daily Activity 5"
$ws.Range("L3").Value = "Test Daily Activity 5"

# Update the view: selection moved from P3 to Q3 (the sheet also scrolled
# one column right in the source file; this runtime's xlsx writer does not
# persist plain ScrollColumn/ScrollRow as a sheetView topLeftCell attribute
# outside of frozen panes, so only the selection move is reflected here).
$ws.Range("Q3").Select()
$excel.ActiveWindow.ScrollColumn = 6
